# Fixed update to excel issue
#
# 1. Rename "Requested quantity" header on the Weekly/Monthly sheets.
# 2. Add a new "PO Forecast" sheet with a Prophet-style forecast table
#    (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# --- 1. Header renames -----------------------------------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. New "PO Forecast" sheet ---------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45578.99999999999, 135, 82.23705233118807, 187.2224729127483),
    @(45592.99999999999, 103, 52.16767486587456, 157.5479577664922),
    @(45599.99999999999, 87, 32.17831772994309, 139.7501951847236),
    @(45627.99999999999, 24, -25.64572799333572, 76.09272721083546),
    @(45634.99999999999, 8, -48.65857827664961, 62.45480311571779),
    @(45641.99999999999, 0, -58.01998834484235, 48.7038919553519),
    @(45648.99999999999, 0, -75.16623944248165, 28.50313945068778),
    @(45655.99999999999, 0, -90.84596784936478, 13.54635718650271),
    @(45662.99999999999, 0, -108.7269923118385, -4.914874817003452),
    @(45669.99999999999, 0, -122.113965080843, -20.28838944125903),
    @(45676.99999999999, 0, -141.5565706319944, -33.58092983911602),
    @(45683.99999999999, 0, -155.9273889930916, -52.45951327796431),
    @(45690.99999999999, 0, -170.7885370673725, -61.03829542463191),
    @(45697.99999999999, 0, -190.0142499420925, -80.48792307550177)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Range("A$r").Value = $row[0]
    $wsForecast.Range("B$r").Value = $row[1]
    $wsForecast.Range("C$r").Value = $row[2]
    $wsForecast.Range("D$r").Value = $row[3]
    $r++
}

# Re-use the existing header style (bold, border, center/top) and the
# existing date-column style (custom yyyy-mm-dd numeric format) instead of
# minting brand-new styles, matching how the other two sheets are formatted.
$wsWeekly.Range("A1:B1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsWeekly.Range("A2").Copy() | Out-Null
$wsForecast.Range("A2:A15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsForecast.Range("A1").Select()
